$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.063.97"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Value = "1.853.37"
$ws.Range("E3").Value = "  +2.49%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.44"
$ws.Range("E5").Value = "  +3.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.622"
$ws.Range("E6").Value = "  +1.42%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.20"
$ws.Range("E8").Value = "  +5.37%  "
$ws.Range("E9").Value = "  +1.80%  "
$ws.Range("E10").Value = "  +1.96%  "
$ws.Range("E11").Value = "  +0.43%  "
$ws.Range("D12").Value = "2.121.33"
$ws.Range("E12").Value = "  +2.55%  "
$ws.Range("D13").Value = "1.881.51"
$ws.Range("E13").Value = "  +3.85%  "
$ws.Range("E14").Value = "  +1.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.678"
$ws.Range("E15").Value = "  +2.14%  "
$ws.Range("E16").Value = "  +2.74%  "
$ws.Range("D17").Value = "35.030.46"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.28"
$ws.Range("D19").Value = "0.0₃0796"
$ws.Range("E19").Value = "  +1.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "240.55"
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.19"
$ws.Range("E21").Value = "  +2.63%  "
$ws.Range("E22").Value = "  +2.28%  "
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("E24").Value = "  +1.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.72"
$ws.Range("E25").Value = "  -1.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.90"
$ws.Range("E26").Value = "  +27.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.98"
$ws.Range("E27").Value = "  +2.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.66"
$ws.Range("E28").Value = "  +2.06%  "
$ws.Range("E29").Value = "  +1.51%  "
$ws.Range("E30").Value = "  +0.22%  "
$ws.Range("E31").Value = "  +2.09%  "
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.04"
$ws.Range("E33").Value = "  +2.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.65"
$ws.Range("E34").Value = "  +23.56%  "
$ws.Range("E35").Value = "  +12.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.32"
$ws.Range("E36").Value = "  +6.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.780"
$ws.Range("E37").Value = "  +12.85%  "
$ws.Range("E38").Value = "  +11.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "91.54"
$ws.Range("E39").Value = "  +1.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0202"
$ws.Range("E40").Value = "  +5.91%  "
$ws.Range("D41").Value = "1.351.62"
$ws.Range("E41").Value = "  +1.43%  "
$ws.Range("E42").Value = "  +5.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "14.77"
$ws.Range("E43").Value = "  +2.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.56"
$ws.Range("E44").Value = "  +52.19%  "
$ws.Range("E45").Value = "  -0.31%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0549"
$ws.Range("E47").Value = "  +7.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.51"
$ws.Range("E48").Value = "  +6.32%  "
$ws.Range("D49").Value = "2.034.60"
$ws.Range("E49").Value = "  +2.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0679"
$ws.Range("E50").Value = "  +2.04%  "
$ws.Range("E51").Value = "  +14.52%  "
